$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 150
$ws.Range("M6").Value = -38

$ws.Range("H12").Value = 225.21428
$ws.Range("I12").Value = 239
$ws.Range("K12").Value = 239
$ws.Range("M12").Value = -69

$ws.Range("H51").Value = 5509.154
$ws.Range("I51").Value = 1812.125
$ws.Range("J51").Value = 7152.278
$ws.Range("K51").Value = 1812.125
$ws.Range("L51").Value = 7152.278
$ws.Range("M51").Value = -1328.125
$ws.Range("N51").Value = -8120.278

$ws.Range("H70").Value = 3517.4285
$ws.Range("J70").Value = 2863.4546
$ws.Range("L70").Value = 8590.363799999999
$ws.Range("N70").Value = -9130.363799999999

$ws.Range("H73").Value = 3517.4285
$ws.Range("J73").Value = 2863.4546
$ws.Range("L73").Value = 8590.363799999999
$ws.Range("N73").Value = -10462.3638

$ws.Range("H80").Value = 2719.3
$ws.Range("I80").Value = 1676.5
$ws.Range("J80").Value = 3166.2144
$ws.Range("K80").Value = 5029.5
$ws.Range("L80").Value = 9498.643199999999
$ws.Range("M80").Value = -4031.5
$ws.Range("N80").Value = -11494.6432

$ws.Range("H83").Value = 2719.3
$ws.Range("I83").Value = 1676.5
$ws.Range("J83").Value = 3166.2144
$ws.Range("K83").Value = 15088.5
$ws.Range("L83").Value = 28495.9296
$ws.Range("M83").Value = -10096.5
$ws.Range("N83").Value = -38479.9296

$ws.Range("H92").Value = 726.82355
$ws.Range("I92").Value = 741
$ws.Range("K92").Value = 741
$ws.Range("M92").Value = 507

$ws.Range("H111").Value = 4500
$ws.Range("I111").Value = 4750
$ws.Range("J111").Value = 4000
$ws.Range("K111").Value = 14250
$ws.Range("L111").Value = 12000
$ws.Range("M111").Value = -11183
$ws.Range("N111").Value = -18134

$ws.Range("H133").Value = 119993.664
$ws.Range("J133").Value = 119993.664
$ws.Range("L133").Value = 119993.664
$ws.Range("N133").Value = -130113.664

$ws.Range("H134").Value = 110000
$ws.Range("J134").Value = 110000
$ws.Range("L134").Value = 110000
$ws.Range("N134").Value = -120140

$ws.Range("H138").Value = 26320378
$ws.Range("J138").Value = 62509650
$ws.Range("L138").Value = 187528950
$ws.Range("N138").Value = -187539230

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3072.6667
$ws.Range("I45").Value = 2166.8333
$ws.Range("K45").Value = 2166.8333
$ws.Range("M45").Value = -1789.8333

$ws.Range("H61").Value = 2629
$ws.Range("I61").Value = 2645.2083
$ws.Range("K61").Value = 2645.2083
$ws.Range("M61").Value = -2433.2083

$ws.Range("H97").Value = 1392.5652
$ws.Range("I97").Value = 1104.4445
$ws.Range("K97").Value = 1104.4445
$ws.Range("M97").Value = -608.4445000000001

$ws.Range("H110").Value = 27662.309
$ws.Range("I110").Value = 27662.309
$ws.Range("K110").Value = 27662.309
$ws.Range("M110").Value = -25617.309

$ws.Range("H136").Value = 2629
$ws.Range("I136").Value = 2645.2083
$ws.Range("K136").Value = 7935.624899999999
$ws.Range("M136").Value = -5385.624899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1129.6333
$ws.Range("I20").Value = 1022.4583
$ws.Range("K20").Value = 1022.4583
$ws.Range("M20").Value = -775.4583

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H94").Value = 1463.0938
$ws.Range("I94").Value = 778.7895
$ws.Range("J94").Value = 2463.2307
$ws.Range("K94").Value = 778.7895
$ws.Range("L94").Value = 2463.2307
$ws.Range("M94").Value = -327.7895
$ws.Range("N94").Value = -3365.2307

$ws.Range("H134").Value = 2885.0833
$ws.Range("I134").Value = 842.44446
$ws.Range("K134").Value = 2527.33338
$ws.Range("M134").Value = 7.666619999999966

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2019.3334
$ws.Range("I16").Value = 1823.2
$ws.Range("K16").Value = 1823.2
$ws.Range("M16").Value = -1536.2

$ws.Range("H113").Value = 2019.3334
$ws.Range("I113").Value = 1823.2
$ws.Range("K113").Value = 1823.2
$ws.Range("M113").Value = 346.8

$ws.Range("H132").Value = 48234.41
$ws.Range("I132").Value = 52552.85
$ws.Range("K132").Value = 157658.55
$ws.Range("M132").Value = -155128.55

$ws.Range("H133").Value = 63890.77
$ws.Range("I133").Value = 33000
$ws.Range("J133").Value = 66465
$ws.Range("K133").Value = 33000
$ws.Range("L133").Value = 66465
$ws.Range("N133").Value = -71525
$ws.Range("M133").Value = -30470

$ws.Range("H134").Value = 1307.5238
$ws.Range("I134").Value = 1097.9
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 3293.7
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -758.7000000000003
$ws.Range("N134").Value = -21570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 892.6667
$ws.Range("I114").Value = 639
$ws.Range("J114").Value = 1400
$ws.Range("K114").Value = 1917
$ws.Range("L114").Value = 4200
$ws.Range("M114").Value = 1337
$ws.Range("N114").Value = -10708

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H80").Value = 2670.5715
$ws.Range("I80").Value = 2159
$ws.Range("K80").Value = 2159
$ws.Range("M80").Value = -1161

$ws.Range("H83").Value = 2670.5715
$ws.Range("I83").Value = 2159
$ws.Range("K83").Value = 10795
$ws.Range("M83").Value = -5803

$ws.Range("H97").Value = 1316.05
$ws.Range("I97").Value = 991.1818
$ws.Range("J97").Value = 1713.1111
$ws.Range("K97").Value = 991.1818
$ws.Range("L97").Value = 1713.1111
$ws.Range("M97").Value = -495.1818
$ws.Range("N97").Value = -2705.1111

$ws.Range("H102").Value = 4899.4346
$ws.Range("I102").Value = 2546.6667
$ws.Range("K102").Value = 2546.6667
$ws.Range("M102").Value = -924.6667000000002

$ws.Range("H122").Value = 921.03125
$ws.Range("I122").Value = 773.3226
$ws.Range("K122").Value = 2319.9678
$ws.Range("M122").Value = 130.0322000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5667.769
$ws.Range("I7").Value = 5281.2856
$ws.Range("J7").Value = 6118.6665
$ws.Range("K7").Value = 5281.2856
$ws.Range("L7").Value = 6118.6665
$ws.Range("M7").Value = -5169.2856
$ws.Range("N7").Value = -6342.6665

$ws.Range("H16").Value = 1760.2
$ws.Range("I16").Value = 1256.7273
$ws.Range("K16").Value = 1256.7273
$ws.Range("M16").Value = -1086.7273

$ws.Range("H40").Value = 1965.6364
$ws.Range("I40").Value = 1662.2
$ws.Range("K40").Value = 1662.2
$ws.Range("M40").Value = -1526.2

$ws.Range("H55").Value = 664.87177
$ws.Range("J55").Value = 822.1
$ws.Range("L55").Value = 822.1
$ws.Range("N55").Value = -1168.1

$ws.Range("H126").Value = 5667.769
$ws.Range("I126").Value = 5281.2856
$ws.Range("J126").Value = 6118.6665
$ws.Range("K126").Value = 15843.8568
$ws.Range("L126").Value = 18355.9995
$ws.Range("M126").Value = -13373.8568
$ws.Range("N126").Value = -23295.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 25000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H81").Value = 3637.4138
$ws.Range("J81").Value = 5399.8
$ws.Range("L81").Value = 10799.6
$ws.Range("N81").Value = -12921.6

$ws.Range("H84").Value = 3637.4138
$ws.Range("J84").Value = 5399.8
$ws.Range("L84").Value = 53998
$ws.Range("N84").Value = -64606

$ws.Range("H107").Value = 1367.7693
$ws.Range("I107").Value = 1516.8334
$ws.Range("K107").Value = 4550.5002
$ws.Range("M107").Value = -2630.5002

$ws.Range("H122").Value = 45719.957
$ws.Range("I122").Value = 55604.684
$ws.Range("J122").Value = 8158
$ws.Range("K122").Value = 166814.052
$ws.Range("L122").Value = 24474
$ws.Range("M122").Value = -164364.052
$ws.Range("N122").Value = -29374
